$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 27: new B27/C27 date entries (mirrors format of the surrounding date cells) ---
$ws.Range("B27").Value = 45513
$ws.Range("C27").Value = 45513
$ws.Range("B2").Copy()
$ws.Range("B27:C27").PasteSpecial(-4122)  # xlPasteFormats - reuse the existing date style

# --- Row 32: DensityGradientDepth use-factors ---
$ws.Range("P32").Formula = "=0.0254*0.0254*0.0254*100*0.3048/0.45359"
$ws.Range("Q32").Formula = "=P32*30/100"
$ws.Range("R32").Formula = "=Q32/30"

# --- Rows 34-37: new unit-test rows (Image Scale / Pressure / Heat Transfer Coefficient / Pressure Loss Constant) ---
$ws.Range("P34").Formula = "=3*0.3048"

$ws.Range("P35").Formula = "=P34*P34*P34*100*0.3048/0.45359"
$ws.Range("Q35").Formula = "=P35*30/100"
$ws.Range("R35").Formula = "=Q35/30"

$ws.Range("P36").Formula = "=0.00454609*100*0.3048/0.45359"
$ws.Range("Q36").Formula = "=P36*30/100"
$ws.Range("R36").Formula = "=Q36/30"

$ws.Range("P37").Formula = "=231*0.0254*0.0254*0.0254*100*0.3048/0.45359"
$ws.Range("Q37").Formula = "=P37*30/100"
$ws.Range("R37").Formula = "=Q37/30"

# --- New column R width (visible now that it holds data) ---
$ws.Columns("R").ColumnWidth = 11.16666666666667

# --- View state: scroll/selection moved as the author kept editing further down the sheet ---
$ws.Range("J28").Select()

$wb.Save()
